$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some "Price" column values are numeric-looking text (e.g. thousand-dot
# separated, or with significant trailing zeros) that must remain stored as
# text exactly as given, so force a Text number format on each such cell
# before assigning the value, then restore the default (unstyled) appearance.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '70.175.48'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -3.03%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.523.06'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -4.24%  '
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '578.53'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.10%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '168.56'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -3.72%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.519'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.01%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.522.22'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -4.20%  '
$ws.Range("E10").Value = '  -5.33%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.168'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.64%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.348'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -2.94%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.92'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.21%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.983.26'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -4.40%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '70.057.40'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -2.96%  '
$ws.Range("E16").Value = '  -5.29%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '25.16'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -2.15%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.531.80'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -4.90%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.79'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.71%  '
$ws.Range("E20").Value = '  -5.71%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '351.50'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -6.55%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.95'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -3.83%  '
$ws.Range("E23").Value = '  -3.94%  '
$ws.Range("E24").Value = '  +0.10%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '69.23'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -3.17%  '
$ws.Range("E26").Value = '  -4.91%  '
$ws.Range("E27").Value = '  -5.01%  '
$ws.Range("E28").Value = '  -4.56%  '
$ws.Range("E29").Value = '  +0.26%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0₃0910'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -3.92%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.93'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.28%  '
$ws.Range("E32").Value = '  -2.03%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '466.00'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -5.12%  '
$ws.Range("E34").Value = '  -1.86%  '
$ws.Range("E35").Value = '  -0.07%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.120'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +4.24%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '153.02'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -5.35%  '
$ws.Range("E38").Value = '  +0.71%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '18.49'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -3.47%  '
$ws.Range("E40").Value = '  +0.02%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '4.80'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.72%  '
$ws.Range("E42").Value = '  -1.25%  '
$ws.Range("E43").Value = '  -6.72%  '
$ws.Range("E44").Value = '  -13.54%  '
$ws.Range("E45").Value = '  -9.35%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '38.22'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.06%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '143.57'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -4.41%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.534'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.37%  '
$ws.Range("E49").Value = '  -3.09%  '
$ws.Range("E50").Value = '  -4.17%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0736'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.17%  '
